# "Common: Liquid mix works quite well... I hope"
#
# Adds two new translation rows to the "Translations - Lab" sheet for the
# new lab.liquid.quick-info.* strings, and makes that sheet the active tab
# with B97 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Lab")

# Insert two blank rows before the current row 92 (pushes old rows 92-110
# down to 94-112, carrying their formatting/styles along for the ride).
$ws.Rows("92:93").Insert()

# Row 93 is filled in first (matches the order new strings were appended
# to the shared-string table), then row 92.
$ws.Cells.Item(93, 1).Value = "cs"
$ws.Cells.Item(93, 2).Value = "lab.liquid.quick-info.title"
$ws.Cells.Item(93, 3).Value = "Liquid"

$ws.Cells.Item(92, 1).Value = "cs"
$ws.Cells.Item(92, 2).Value = "lab.liquid.quick-info.subtitle"
$ws.Cells.Item(92, 3).Value = "Zde se zobrazí informace o složení liquidu."

# Make "Translations - Lab" the active sheet/tab and select B97.
$ws.Activate()
$ws.Range("B97").Select()
